$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.100939273834229
$ws.Range("B1").Value = 2.69025444984436
$ws.Range("C1").Value = 2.817869663238525
$ws.Range("D1").Value = 2.912649393081665
$ws.Range("E1").Value = 0.7650070786476135
